# Rename worksheets
$wb = $excel.ActiveWorkbook

$wsBase = $wb.Worksheets.Item(1)
$wsBase.Name = "Base"

$wsExpansion = $wb.Worksheets.Item(2)
$wsExpansion.Name = "Expansion"

$wsResources = $wb.Worksheets.Item(3)
$wsResources.Name = "Resources"

# Populate Expansion sheet
$wsExpansion.Range("A1").Value = "Name"
$wsExpansion.Range("B1").Value = "Level"
$wsExpansion.Range("A1:B1").Font.Bold = $true

$wsExpansion.Range("A2").Value = "Bootlegger"
$wsExpansion.Range("B2").Value = 1
$wsExpansion.Range("A3").Value = "Rum Runner"
$wsExpansion.Range("B3").Value = 2
$wsExpansion.Range("A4").Value = "Whiskey Baron"
$wsExpansion.Range("B4").Value = 3

# Populate Resources sheet (write column A fully before column B so that
# shared-string indices are allocated in the same order as the target file:
# Wood, Metal, Stone, then Cost)
$wsResources.Range("A1").Value = "Name"
$wsResources.Range("A2").Value = "Wood"
$wsResources.Range("A3").Value = "Metal"
$wsResources.Range("A4").Value = "Stone"

$wsResources.Range("B1").Value = "Cost"
$wsResources.Range("B2").Value = 2
$wsResources.Range("B3").Value = 3
$wsResources.Range("B4").Value = 5

$wsResources.Range("A1:B1").Font.Bold = $true

# Resources sheet gains a page setup (portrait orientation), mirroring sheet1's existing pageSetup
$wsResources.PageSetup.Orientation = 1

# Selections: clear selection on Base, select B4 on Expansion, make Resources the active (selected) tab
$wsBase.Range("A1:C4").Select()
$wsExpansion.Range("B4").Select()

# Make Resources the active sheet (tab selected)
$wsResources.Activate()
